$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are stored as text, matching the source data
# (values such as "1.005" or "6.060" would otherwise be auto-converted to numbers)
$ws.Range('D2:D51').NumberFormat = '@'

# Row 2
$ws.Range('D2').Value = '25.892.83'
$ws.Range('E2').Value = '  +0.20%  '
# Row 3
$ws.Range('D3').Value = '1.644.87'
$ws.Range('E3').Value = '  +0.62%  '
# Row 4
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.35%  '
# Row 5
$ws.Range('D5').Value = '215.29'
# Row 6
$ws.Range('D6').Value = '0.5073'
$ws.Range('E6').Value = '  +0.84%  '
# Row 7
$ws.Range('D7').Value = '1.005'
$ws.Range('E7').Value = '  +0.30%  '
# Row 8
$ws.Range('D8').Value = '0.2574'
$ws.Range('E8').Value = '  +0.01%  '
# Row 9
$ws.Range('D9').Value = '0.06411'
$ws.Range('E9').Value = '  +0.32%  '
# Row 10
$ws.Range('D10').Value = '19.75'
$ws.Range('E10').Value = '  +0.53%  '
# Row 11
$ws.Range('D11').Value = '0.07768'
$ws.Range('E11').Value = '  +1.30%  '
# Row 12
$ws.Range('D12').Value = '4.307'
$ws.Range('E12').Value = '  +1.55%  '
# Row 13
$ws.Range('D13').Value = '1.651.73'
$ws.Range('E13').Value = '  +0.66%  '
# Row 14
$ws.Range('D14').Value = '0.5464'
$ws.Range('E14').Value = '  +0.22%  '
# Row 15
$ws.Range('D15').Value = '0.0₅7893'
$ws.Range('E15').Value = '  -0.43%  '
# Row 16
$ws.Range('E16').Value = '  +2.48%  '
# Row 17
$ws.Range('D17').Value = '25.977.72'
$ws.Range('E17').Value = '  +0.50%  '
# Row 18
$ws.Range('D18').Value = '1.006'
$ws.Range('E18').Value = '  +0.39%  '
# Row 19
$ws.Range('D19').Value = '197.55'
$ws.Range('E19').Value = '  -2.65%  '
# Row 20
$ws.Range('D20').Value = '4.412'
$ws.Range('E20').Value = '  +2.11%  '
# Row 21
$ws.Range('D21').Value = '10.01'
$ws.Range('E21').Value = '  +0.73%  '
# Row 22
$ws.Range('D22').Value = '6.060'
$ws.Range('E22').Value = '  +1.60%  '
# Row 23
$ws.Range('E23').Value = '  +0.50%  '
# Row 24
$ws.Range('D24').Value = '1.860'
$ws.Range('E24').Value = '  -3.98%  '
# Row 25
$ws.Range('D25').Value = '141.07'
$ws.Range('E25').Value = '  +0.16%  '
# Row 26
$ws.Range('E26').Value = '  +0.32%  '
# Row 27
$ws.Range('E27').Value = '  +2.97%  '
# Row 28
$ws.Range('D28').Value = '15.73'
$ws.Range('E28').Value = '  +0.18%  '
# Row 29
$ws.Range('D29').Value = '0.05082'
$ws.Range('E29').Value = '  +2.07%  '
# Row 30
$ws.Range('D30').Value = '1.241'
$ws.Range('E30').Value = '  +0.10%  '
# Row 31
$ws.Range('D31').Value = '3.270'
$ws.Range('E31').Value = '  -0.16%  '
# Row 32
$ws.Range('D32').Value = '3.206'
$ws.Range('E32').Value = '  +0.61%  '
# Row 33
$ws.Range('D33').Value = '1.541'
$ws.Range('E33').Value = '  +0.28%  '
# Row 34
$ws.Range('D34').Value = '2.366'
$ws.Range('E34').Value = '  +0.64%  '
# Row 35
$ws.Range('D35').Value = '0.8942'
$ws.Range('E35').Value = '  +0.14%  '
# Row 36
$ws.Range('D36').Value = '2.596'
$ws.Range('E36').Value = '  -0.85%  '
# Row 37
$ws.Range('D37').Value = '0.5549'
$ws.Range('E37').Value = '  -0.62%  '
# Row 38
$ws.Range('D38').Value = '1.134.94'
$ws.Range('E38').Value = '  -3.53%  '
# Row 39
$ws.Range('D39').Value = '0.01565'
$ws.Range('E39').Value = '  +0.48%  '
# Row 40
$ws.Range('D40').Value = '1.006'
$ws.Range('E40').Value = '  +0.44%  '
# Row 41
$ws.Range('E41').Value = '  +0.43%  '
# Row 42
$ws.Range('D42').Value = '0.8152'
$ws.Range('E42').Value = '  +1.20%  '
# Row 43
$ws.Range('D43').Value = '99.72'
$ws.Range('E43').Value = '  +0.41%  '
# Row 44
$ws.Range('D44').Value = '0.0₈121'
$ws.Range('E44').Value = '  +6.03%  '
# Row 45
$ws.Range('D45').Value = '1.782.27'
$ws.Range('E45').Value = '  +0.66%  '
# Row 46
$ws.Range('D46').Value = '0.4530'
$ws.Range('E46').Value = '  +0.40%  '
# Row 47
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').Value = '1.008'
$ws.Range('E47').Value = '  +0.17%  '
# Row 48
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '55.27'
$ws.Range('E48').Value = '  +0.76%  '
# Row 49
$ws.Range('D49').Value = '0.05082'
$ws.Range('E49').Value = '  +1.01%  '
# Row 50
$ws.Range('D50').Value = '1.007'
$ws.Range('E50').Value = '  +0.47%  '
# Row 51
$ws.Range('D51').Value = '0.09571'
$ws.Range('E51').Value = '  +3.09%  '
